$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stable donor cells we can copy formatting (style index) from without
# disturbing them: B13 keeps style 2 ("Menlo" xf) for the whole edit, and
# G14 keeps style 1 ("Arial" xf) for the whole edit.
$style2Donor = $ws.Range("B13")
$style1Donor = $ws.Range("G14")

function Set-StyledValue($cellRef, $donor, $value) {
    $donor.Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
    $ws.Range($cellRef).Value = $value
}

# --- Header row 11 -------------------------------------------------------
# t1/t2/t3/t-avg/Fish (H:L) collapses down to just "parallel" (H) and
# "sequential" (I); style 1 for both; J:L removed entirely.
Set-StyledValue "H11" $style1Donor "parallel"
Set-StyledValue "I11" $style1Donor "sequential"
$ws.Range("J11:L11").Clear()

# --- Row 12 (100 fish) -----------------------------------------------------
Set-StyledValue "H12" $style1Donor 1.4437
Set-StyledValue "I12" $style2Donor 0.0494
$ws.Range("J12:L12").Clear()

# --- Row 13 (1000 fish) -----------------------------------------------------
Set-StyledValue "H13" $style1Donor 3.6846
Set-StyledValue "I13" $style1Donor 0.2472
$ws.Range("J13:L13").Clear()

# --- Row 14 (10000 fish) -----------------------------------------------------
Set-StyledValue "H14" $style1Donor 23.541971
Set-StyledValue "I14" $style1Donor 1.8972
$ws.Range("J14:L14").Clear()

# --- Row 15 (100000 fish) -----------------------------------------------------
Set-StyledValue "H15" $style1Donor 209.4251
Set-StyledValue "I15" $style2Donor 18.4586
$ws.Range("J15:L15").Clear()

# --- Row 16 (1000000 fish) -----------------------------------------------------
Set-StyledValue "H16" $style1Donor 1574.5107
Set-StyledValue "I16" $style2Donor 181.5163
$ws.Range("J16:L16").Clear()

# --- Row 17 (1500000 fish) -----------------------------------------------------
Set-StyledValue "H17" $style2Donor 1964.2471
Set-StyledValue "I17" $style2Donor 273.5627
$ws.Range("J17:L17").Clear()
# J17 keeps an (empty) styled cell in the target workbook.
$style2Donor.Copy()
$ws.Range("J17").PasteSpecial(-4122)

$ws.Range("N11").Select()
